# Updated cryptos list on Sat Jun  8 08:49:06 UTC 2024 with GitHub Actions
#
# The "Price" column (D) stores numbers as plain text (e.g. "69.377.42",
# "693.27") because some values use "." as a thousands separator, which is
# not valid numeric syntax. Assigning a plain numeric-looking string such
# as "693.27" via .Value would normally get auto-coerced to a real number
# by Excel, which would break the text formatting used throughout this
# column. To avoid that, the whole Price column is temporarily switched to
# a Text number format before the new values are written, then the
# temporary formatting is cleared again so the cells end up unstyled, same
# as before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Row 2 : Bitcoin ---
$ws.Range("D2").Value = "69.377.42"
$ws.Range("E2").Value = "  -2.39%  "

# --- Row 3 : Ethereum ---
$ws.Range("D3").Value = "3.693.43"
$ws.Range("E3").Value = "  -3.09%  "

# --- Row 4 : TetherUSD ---
$ws.Range("E4").Value = "  +0.05%  "

# --- Row 5 : BNB ---
$ws.Range("D5").Value = "693.27"
$ws.Range("E5").Value = "  -1.11%  "

# --- Row 6 : Solana ---
$ws.Range("D6").Value = "162.63"
$ws.Range("E6").Value = "  -5.56%  "

# --- Row 7 : LidoStakedEther ---
$ws.Range("D7").Value = "3.693.01"
$ws.Range("E7").Value = "  -3.08%  "

# --- Row 8 : USDC ---
$ws.Range("E8").Value = "  +0.02%  "

# --- Row 9 : XRP ---
$ws.Range("E9").Value = "  -4.66%  "

# --- Row 11 : Toncoin ---
$ws.Range("E11").Value = "  -1.97%  "

# --- Row 12 : Cardano ---
$ws.Range("E12").Value = "  -4.58%  "

# --- Row 13 : ShibaInu ---
$ws.Range("E13").Value = "  -5.23%  "

# --- Row 14 : Avalanche ---
$ws.Range("D14").Value = "33.39"
$ws.Range("E14").Value = "  -7.32%  "

# --- Row 15 : WrappedliquidstakedEther2.0 ---
$ws.Range("D15").Value = "4.314.50"
$ws.Range("E15").Value = "  -3.12%  "

# --- Row 16 : WrappedEther ---
$ws.Range("D16").Value = "3.690.63"
$ws.Range("E16").Value = "  -3.17%  "

# --- Row 17 : WrappedBTC ---
$ws.Range("D17").Value = "69.417.03"
$ws.Range("E17").Value = "  -2.37%  "

# --- Row 18 : TRON ---
$ws.Range("E18").Value = "  -0.89%  "

# --- Row 19 : Chainlink ---
$ws.Range("D19").Value = "16.24"
$ws.Range("E19").Value = "  -7.26%  "

# --- Row 20 : Polkadot ---
$ws.Range("D20").Value = "6.59"
$ws.Range("E20").Value = "  -7.73%  "

# --- Row 21 : BitcoinCash ---
$ws.Range("D21").Value = "480.06"
$ws.Range("E21").Value = "  -6.53%  "

# --- Row 22 : Uniswap ---
$ws.Range("E22").Value = "  -6.33%  "

# --- Row 23 : Polygon ---
$ws.Range("D23").Value = "0.663"
$ws.Range("E23").Value = "  -7.32%  "

# --- Row 24 : Litecoin ---
$ws.Range("D24").Value = "80.04"
$ws.Range("E24").Value = "  -4.69%  "

# --- Row 25 : WrappedeETH ---
$ws.Range("D25").Value = "3.839.84"
$ws.Range("E25").Value = "  -3.12%  "

# --- Row 26 : PEPE ---
$ws.Range("D26").Value = "0.0000130"
$ws.Range("E26").Value = "  -9.66%  "

# --- Row 27 : Dai ---
$ws.Range("E27").Value = "  +0.04%  "

# --- Row 28 : InternetComputer(DFINITY) ---
$ws.Range("D28").Value = "11.36"
$ws.Range("E28").Value = "  -5.99%  "

# --- Row 29 : RenderToken ---
$ws.Range("E29").Value = "  -8.79%  "

# --- Row 30 : Fetch.AI ---
$ws.Range("E30").Value = "  -10.99%  "

# --- Row 31 : PancakeSwap ---
$ws.Range("D31").Value = "2.72"
$ws.Range("E31").Value = "  -10.15%  "

# --- Row 32 : NEARProtocol ---
$ws.Range("D32").Value = "6.83"
$ws.Range("E32").Value = "  -8.06%  "

# --- Row 33 : ImmutableX ---
$ws.Range("E33").Value = "  -7.75%  "

# --- Rows 34-36: coins reshuffled (Kaspa -> EthereumClassic -> Binance-PegBSC-USD -> Kaspa) ---
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "27.00"
$ws.Range("E34").Value = "  -7.12%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.167"
$ws.Range("E36").Value = "  -4.40%  "

# --- Row 37 : RenzoRestakedETH ---
$ws.Range("D37").Value = "3.661.07"
$ws.Range("E37").Value = "  -2.95%  "

# --- Row 38 : Aptos ---
$ws.Range("D38").Value = "8.48"
$ws.Range("E38").Value = "  -7.54%  "

# --- Row 39 : Filecoin ---
$ws.Range("D39").Value = "6.37"
$ws.Range("E39").Value = "  +5.80%  "

# --- Row 40 : Stacks ---
$ws.Range("D40").Value = "2.33"
$ws.Range("E40").Value = "  -2.13%  "

# --- Row 41 : Hedera ---
$ws.Range("D41").Value = "0.0930"
$ws.Range("E41").Value = "  -7.98%  "

# --- Row 42 : USDe ---
$ws.Range("E42").Value = "  +0.00%  "

# --- Row 43 : FirstDigitalUSD ---
$ws.Range("E43").Value = "  -0.04%  "

# --- Row 44 : Mantle ---
$ws.Range("D44").Value = "0.956"
$ws.Range("E44").Value = "  -6.39%  "

# --- Row 45 : Monero ---
$ws.Range("D45").Value = "164.11"
$ws.Range("E45").Value = "  -5.14%  "

# --- Row 46 : OKB ---
$ws.Range("D46").Value = "47.98"
$ws.Range("E46").Value = "  -2.84%  "

# --- Row 47 : InjectiveProtocol ---
$ws.Range("D47").Value = "30.17"
$ws.Range("E47").Value = "  +2.62%  "

# --- Row 48 : dogwifhat ---
$ws.Range("E48").Value = "  -15.28%  "

# --- Row 49 : ONDO ---
$ws.Range("E49").Value = "  -1.51%  "

# --- Rows 50-51: coins swapped (FLOKI <-> SuiNetwork) ---
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "1.14"
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "0.000285"
$ws.Range("E51").Value = "  -9.13%  "

# Restore the Price column cells to their original (unstyled) state now
# that the text values are safely stored.
$priceRange.ClearFormats()
